# Auto-applied edit based on diff: updates Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks like a plain number (e.g. "210.12") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers and
# drops significant trailing zeros (e.g. "0.890" -> 0.89).
$textForceCells = @(
    'D5',
    'D6',
    'D8',
    'D11',
    'D15',
    'D17',
    'D18',
    'D25',
    'D40',
    'D43',
    'D44',
    'D46',
    'D48',
    'D49',
    'D50',
)
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.325.62'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '1.552.51'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '210.12'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = '0.481'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '23.85'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.774.95'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').Value = '1.548.65'
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').Value = '28.326.28'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').Value = '0.511'
$ws.Range('E15').Value = '  -2.07%  '
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = '60.85'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').Value = '227.99'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = '150.91'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('E29').Value = '  -3.23%  '
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('E31').Value = '  -4.64%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '1.388.22'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  -2.83%  '
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').Value = '0.514'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').Value = '0.776'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('D44').Value = '0.0459'
$ws.Range('E44').Value = '  -1.18%  '
$ws.Range('E45').Value = '  -2.18%  '
$ws.Range('D46').Value = '61.94'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').Value = '1.687.16'
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('D48').Value = '0.890'
$ws.Range('E48').Value = '  -7.25%  '
$ws.Range('D49').Value = '85.61'
$ws.Range('E49').Value = '  -1.05%  '
$ws.Range('D50').Value = '42.59'
$ws.Range('E50').Value = '  +7.33%  '
$ws.Range('E51').Value = '  -0.14%  '
